# Updated cryptos list (GitHub Actions data refresh).
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns on
# the cryptos sheet with newly scraped values. All D/E cells in this sheet
# are stored as plain text (coinranking.com formats prices/volumes as
# strings, not numbers), so every write below keeps the cell as Text.
#
# A handful of the new Price values (column D) look like plain numbers
# (e.g. "0.0846", "210.73") and Excel would otherwise auto-convert them to
# the Number type on assignment. For those we prefix the literal with a
# leading apostrophe (Excel's standard "force text" marker) and then
# reapply the Normal cell style so the stored style index is left exactly
# as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.303.82"
$ws.Range("D3").Value = "1.588.33"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "'0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "1.601.58"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "26.311.69"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("E19").Value = "  +5.95%  "
$ws.Range("D20").Value = "'210.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'144.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D30").Value = "'0.0506"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "1.312.57"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -10.56%  "
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("D43").Value = "'0.767"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'2.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").Value = "'62.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").Value = "1.724.84"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "'87.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("E48").Value = "  -5.37%  "
$ws.Range("D49").Value = "'0.0507"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").Value = "'0.0979"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("E51").Value = "  -0.40%  "
